$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the A1 conversion summary text ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.74 = 51847.13 pesos`n✅ 51847.13 pesos = 12.71 = 927.66 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the N10/O10/N12/O12 rate values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 78.5
$ws2.Range("O10").Value = 4070
$ws2.Range("N12").Value = 4079.98
$ws2.Range("O12").Value = 73
